$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value2 = "FAPs"
$ws.Cells.Item(2,2).Value2 = "Wnt5a"
$ws.Cells.Item(2,3).Value2 = "Fzd6"
$ws.Cells.Item(2,4).Value2 = "ECs"
$ws.Cells.Item(2,5).Value2 = 3
$ws.Cells.Item(2,6).Value2 = 1
$ws.Cells.Item(2,7).Value2 = 7.514794999999999
$ws.Cells.Item(2,8).Value2 = 22.544385
$ws.Cells.Item(2,9).Value2 = 0.992147452492356
$ws.Cells.Item(2,10).Value2 = 0.992147452492356
$ws.Cells.Item(2,11).Value2 = 3
$ws.Cells.Item(2,12).Value2 = 1
$ws.Cells.Item(2,13).Value2 = 13.877148
$ws.Cells.Item(2,14).Value2 = 41.631444
$ws.Cells.Item(2,15).Value2 = 0.9551716529386821
$ws.Cells.Item(2,16).Value2 = 0.9551716529386822
$ws.Cells.Item(2,17).Value2 = 104.28392240466
$ws.Cells.Item(2,18).Value2 = 938.55530164194
$ws.Cells.Item(2,19).Value2 = 0.9476711221560263
$ws.Cells.Item(2,20).Value2 = 0.9476711221560264

# Row 3
$ws.Cells.Item(3,1).Value2 = "FAPs"
$ws.Cells.Item(3,2).Value2 = "Wnt5a"
$ws.Cells.Item(3,3).Value2 = "Fzd6"
$ws.Cells.Item(3,4).Value2 = "FAPs"
$ws.Cells.Item(3,5).Value2 = 3
$ws.Cells.Item(3,6).Value2 = 1
$ws.Cells.Item(3,7).Value2 = 7.514794999999999
$ws.Cells.Item(3,8).Value2 = 22.544385
$ws.Cells.Item(3,9).Value2 = 0.992147452492356
$ws.Cells.Item(3,10).Value2 = 0.992147452492356
$ws.Cells.Item(3,11).Value2 = 2
$ws.Cells.Item(3,12).Value2 = 0.6666666666666666
$ws.Cells.Item(3,13).Value2 = 0.4259926666666667
$ws.Cells.Item(3,14).Value2 = 1.277978
$ws.Cells.Item(3,15).Value2 = 0.02932130719941569
$ws.Cells.Item(3,16).Value2 = 0.02932130719941569
$ws.Cells.Item(3,17).Value2 = 3.201247561503333
$ws.Cells.Item(3,18).Value2 = 28.81122805353
$ws.Cells.Item(3,19).Value2 = 0.02909106024164606
$ws.Cells.Item(3,20).Value2 = 0.02909106024164606

# Row 4
$ws.Cells.Item(4,1).Value2 = "FAPs"
$ws.Cells.Item(4,2).Value2 = "Wnt5a"
$ws.Cells.Item(4,3).Value2 = "Fzd6"
$ws.Cells.Item(4,4).Value2 = "sCs"
$ws.Cells.Item(4,5).Value2 = 3
$ws.Cells.Item(4,6).Value2 = 1
$ws.Cells.Item(4,7).Value2 = 7.514794999999999
$ws.Cells.Item(4,8).Value2 = 22.544385
$ws.Cells.Item(4,9).Value2 = 0.992147452492356
$ws.Cells.Item(4,10).Value2 = 0.992147452492356
$ws.Cells.Item(4,11).Value2 = 3
$ws.Cells.Item(4,12).Value2 = 1
$ws.Cells.Item(4,13).Value2 = 0.225293
$ws.Cells.Item(4,14).Value2 = 0.675879
$ws.Cells.Item(4,15).Value2 = 0.01550703986190206
$ws.Cells.Item(4,16).Value2 = 0.01550703986190207
$ws.Cells.Item(4,17).Value2 = 1.693030709935
$ws.Cells.Item(4,18).Value2 = 15.237276389415
$ws.Cells.Item(4,19).Value2 = 0.01538527009468355
$ws.Cells.Item(4,20).Value2 = 0.01538527009468355

# Row 5
$ws.Cells.Item(5,1).Value2 = "sCs"
$ws.Cells.Item(5,2).Value2 = "Wnt5a"
$ws.Cells.Item(5,3).Value2 = "Fzd6"
$ws.Cells.Item(5,4).Value2 = "ECs"
$ws.Cells.Item(5,5).Value2 = 1
$ws.Cells.Item(5,6).Value2 = 0.3333333333333333
$ws.Cells.Item(5,7).Value2 = 0.05947733333333333
$ws.Cells.Item(5,8).Value2 = 0.178432
$ws.Cells.Item(5,9).Value2 = 0.007852547507643968
$ws.Cells.Item(5,10).Value2 = 0.00785254750764397
$ws.Cells.Item(5,11).Value2 = 3
$ws.Cells.Item(5,12).Value2 = 1
$ws.Cells.Item(5,13).Value2 = 13.877148
$ws.Cells.Item(5,14).Value2 = 41.631444
$ws.Cells.Item(5,15).Value2 = 0.9551716529386821
$ws.Cells.Item(5,16).Value2 = 0.9551716529386822
$ws.Cells.Item(5,17).Value2 = 0.825375757312
$ws.Cells.Item(5,18).Value2 = 7.428381815808001
$ws.Cells.Item(5,19).Value2 = 0.007500530782655818
$ws.Cells.Item(5,20).Value2 = 0.00750053078265582

# Row 6
$ws.Cells.Item(6,1).Value2 = "sCs"
$ws.Cells.Item(6,2).Value2 = "Wnt5a"
$ws.Cells.Item(6,3).Value2 = "Fzd6"
$ws.Cells.Item(6,4).Value2 = "FAPs"
$ws.Cells.Item(6,5).Value2 = 1
$ws.Cells.Item(6,6).Value2 = 0.3333333333333333
$ws.Cells.Item(6,7).Value2 = 0.05947733333333333
$ws.Cells.Item(6,8).Value2 = 0.178432
$ws.Cells.Item(6,9).Value2 = 0.007852547507643968
$ws.Cells.Item(6,10).Value2 = 0.00785254750764397
$ws.Cells.Item(6,11).Value2 = 2
$ws.Cells.Item(6,12).Value2 = 0.6666666666666666
$ws.Cells.Item(6,13).Value2 = 0.4259926666666667
$ws.Cells.Item(6,14).Value2 = 1.277978
$ws.Cells.Item(6,15).Value2 = 0.02932130719941569
$ws.Cells.Item(6,16).Value2 = 0.02932130719941569
$ws.Cells.Item(6,17).Value2 = 0.02533690783288889
$ws.Cells.Item(6,18).Value2 = 0.228032170496
$ws.Cells.Item(6,19).Value2 = 0.0002302469577696348
$ws.Cells.Item(6,20).Value2 = 0.0002302469577696349

# Row 7
$ws.Cells.Item(7,1).Value2 = "sCs"
$ws.Cells.Item(7,2).Value2 = "Wnt5a"
$ws.Cells.Item(7,3).Value2 = "Fzd6"
$ws.Cells.Item(7,4).Value2 = "sCs"
$ws.Cells.Item(7,5).Value2 = 1
$ws.Cells.Item(7,6).Value2 = 0.3333333333333333
$ws.Cells.Item(7,7).Value2 = 0.05947733333333333
$ws.Cells.Item(7,8).Value2 = 0.178432
$ws.Cells.Item(7,9).Value2 = 0.007852547507643968
$ws.Cells.Item(7,10).Value2 = 0.00785254750764397
$ws.Cells.Item(7,11).Value2 = 3
$ws.Cells.Item(7,12).Value2 = 1
$ws.Cells.Item(7,13).Value2 = 0.225293
$ws.Cells.Item(7,14).Value2 = 0.675879
$ws.Cells.Item(7,15).Value2 = 0.01550703986190206
$ws.Cells.Item(7,16).Value2 = 0.01550703986190207
$ws.Cells.Item(7,17).Value2 = 0.01339982685866667
$ws.Cells.Item(7,18).Value2 = 0.120598441728
$ws.Cells.Item(7,19).Value2 = 0.0001217697672185147
$ws.Cells.Item(7,20).Value2 = 0.0001217697672185148
